$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$c = $t.Cell(1, 1)
$c.Range.Text = "45 x 83" + $vt + "  8    3" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"

$c = $t.Cell(1, 2)
$c.Range.Text = "79 x 46" + $vt + "  4    6" + $vt + "  ----" + $vt + "7|    |" + $vt + "9|    |"

$c = $t.Cell(1, 3)
$c.Range.Text = "19 x 54" + $vt + "  5    4" + $vt + "  ----" + $vt + "1|    |" + $vt + "9|    |"

$c = $t.Cell(2, 1)
$c.Range.Text = "32 x 43" + $vt + "  4    3" + $vt + "  ----" + $vt + "3|    |" + $vt + "2|    |"

$c = $t.Cell(2, 2)
$c.Range.Text = "26 x 27" + $vt + "  2    7" + $vt + "  ----" + $vt + "2|    |" + $vt + "6|    |"

$c = $t.Cell(2, 3)
$c.Range.Text = "78 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "7|    |" + $vt + "8|    |"

$c = $t.Cell(3, 1)
$c.Range.Text = "66 x 44" + $vt + "  4    4" + $vt + "  ----" + $vt + "6|    |" + $vt + "6|    |"

$c = $t.Cell(3, 2)
$c.Range.Text = "54 x 93" + $vt + "  9    3" + $vt + "  ----" + $vt + "5|    |" + $vt + "4|    |"

$c = $t.Cell(3, 3)
$c.Range.Text = "70 x 94" + $vt + "  9    4" + $vt + "  ----" + $vt + "7|    |" + $vt + "0|    |"

$c = $t.Cell(4, 1)
$c.Range.Text = "54 x 71" + $vt + "  7    1" + $vt + "  ----" + $vt + "5|    |" + $vt + "4|    |"

$c = $t.Cell(4, 2)
$c.Range.Text = "90 x 46" + $vt + "  4    6" + $vt + "  ----" + $vt + "9|    |" + $vt + "0|    |"

$c = $t.Cell(4, 3)
$c.Range.Text = "62 x 58" + $vt + "  5    8" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"

$c = $t.Cell(5, 1)
$c.Range.Text = "91 x 24" + $vt + "  2    4" + $vt + "  ----" + $vt + "9|    |" + $vt + "1|    |"

$c = $t.Cell(5, 2)
$c.Range.Text = "62 x 82" + $vt + "  8    2" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"

$c = $t.Cell(5, 3)
$c.Range.Text = "24 x 83" + $vt + "  8    3" + $vt + "  ----" + $vt + "2|    |" + $vt + "4|    |"
